# Fix capitalization of the "ProductTb" sheet name -> "ProductTB"
# and make that sheet the active/selected sheet (it was previously "OrderTB").

$wb = $excel.ActiveWorkbook

# Rename sheet to fix capitalization error.
$productSheet = $wb.Worksheets.Item("ProductTb")
$productSheet.Name = "ProductTB"

# Make ProductTB the active sheet, which also moves the "tabSelected" flag
# from the previously active sheet (OrderTB) to this one, and updates the
# workbook's activeTab index.
$productSheet.Activate()
